# Add a new Job Posting row with Job_Id = JD_003
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "JD_003"
$ws.Range("B4").Value = "Senior Engineer"
$ws.Range("C4").Value = "ewfer ewrwerfew"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 4
